$p = $ppt.ActivePresentation
Write-Output $p.Designs.Count
for ($i=1; $i -le $p.Designs.Count; $i++) {
    $d = $p.Designs.Item($i)
    Write-Output "Design $i : Name=$($d.Name) Index=$($d.Index)"
    Write-Output "  SlideMaster.Name=$($d.SlideMaster.Name)"
    Write-Output "  Theme.Name=$($d.SlideMaster.Theme.Name)"
}
